$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds the "last changed" date for every logging
# notice row (2-83). The source feed refreshed and every row's changed-date
# moved forward by 10 days: 2023-11-03 (45233) -> 2023-11-13 (45243).
$ws.Range("C2:C83").Value = 45243
